# Applies the "Make changes to UG/DG and shorten PPP" edit:
#   - shrink the "parseCommand(select)" label font (autofit box shrinks too)
#   - narrow the big green "DG" rounded-rect container
#   - move/resize & re-wrap the ":RestaurantSummaryPanel" label into
#     ":Restaurant" / "SummaryPanel" on two lines
#   - nudge the "loadSummary(restaurant)" label

function Get-ShapeById {
    param($slide, $id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape id=80 "TextBox 79" -> "parseCommand(\u201cselect\u201d)"
# Font size 12 -> 11 on both runs; spAutoFit recalculates the box height
# from 184666 EMU to 169277 EMU.
# ---------------------------------------------------------------------
$sh80 = Get-ShapeById $s 80
$sh80.TextFrame.TextRange.Font.Size = 11
$sh80.Height = 13.328898037795275

# ---------------------------------------------------------------------
# Shape id=59 "Rectangle 65" -> big rounded rectangle container
# Width 3488556 -> 3417931 EMU (position/height unchanged).
# ---------------------------------------------------------------------
$sh59 = Get-ShapeById $s 59
$sh59.Width = 269.1284332968504

# ---------------------------------------------------------------------
# Shape id=71 "Rectangle 62" -> ":RestaurantSummaryPanel" label
# Reposition/resize, and split the text into two centred paragraphs:
#   ":Restaurant" / "SummaryPanel"
# ---------------------------------------------------------------------
$sh71 = Get-ShapeById $s 71
$sh71.Left = 738.7784251968504
$sh71.Top = 250.8427559055118
$sh71.Width = 121.80803299606299
$sh71.Height = 41.77984431968504

$tr71 = $sh71.TextFrame.TextRange
$run1 = $tr71.Runs(1)
[void]$run1.InsertAfter([char]13)

$para1 = $tr71.Paragraphs(1)
$para1Run1 = $para1.Runs(1)
$para1Run2 = $para1.Runs(2)
$para1Run1.Text = ":Restaurant"
$para1Run2.Text = ""

$para2 = $tr71.Paragraphs(2)
$para2Run1 = $para2.Runs(1)
$para2Run1.Text = "SummaryPanel"

# ---------------------------------------------------------------------
# Shape id=54 "TextBox 53" -> "loadSummary(restaurant)" label
# Reposition only (size unchanged).
# ---------------------------------------------------------------------
$sh54 = Get-ShapeById $s 54
$sh54.Left = 720.9704895409449
$sh54.Top = 297.6003265606299
